$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This test-report sheet has one "highlighted status" local format (the fill
# used together with border/alignment) that is currently applied to I11
# ("FAILED"). The edit moves that highlighted look from row 11 (which becomes
# a plain "Not Run" row, like every other untouched row) to row 14 (which
# becomes the "PASSED" row) and, at the same time, re-colors the highlight
# from red (FFC7CE) to green (C6EFCE).
# ---------------------------------------------------------------------------

# 1) Recolor the shared highlight fill while I11 is still its only owner, so
#    the in-place style (currently used only by I11) is the one that picks up
#    the new green color.
#    C6EFCE in BGR (OLE color) order -> 0xCEEFC6
$ws.Range("I11").Interior.Color = 0xCEEFC6

# 2) Propagate that (now green) highlighted format to I14, then set its text.
$ws.Range("I11").Copy($ws.Range("I14"))
$ws.Range("H14").Value = "Reset Pass functionality verified"
$ws.Range("I14").Value = "PASSED"

# 3) Turn I11 back into a plain, un-highlighted "Not Run" cell (matching the
#    formatting used by every other row), and set its text.
$ws.Range("I12").Copy($ws.Range("I11"))
$ws.Range("H11").Value = "Test not executed"
$ws.Range("I11").Value = "Not Run"
